# Steve-1: Fix files BSA-Decomposition
#
# The catalogue rows for the first three information sources still showed the
# older "Актуальность" update stamp (date "BSA 00, 10.01.25" / author
# "Sigfrydj"), while every other row already carries the newer stamp
# ("BSA 00, 13.01.25" / "carolyeu"). Bring rows 3-5 in line with the rest of
# the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: F = Актуальность (unchanged), G = Дата добавления, H = Автор изменений каталога
$ws.Range("G3:G5").Value2 = "BSA 00, 13.01.25"
$ws.Range("H3:H5").Value2 = "carolyeu"

# Refresh the view so the active cell/selection matches the reviewed rows
$ws.Activate()
$ws.Range("I5").Select()
